$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pipe-delimited table of "CellRef|NewValue" pairs describing every cell
# touched by this edit (crypto price/volume refresh + a BitrueCoin /
# MandalaExchangeToken row swap). Columns B/C hold plain text (coin name /
# link) so a normal .Value assignment is safe. Columns D/E hold
# numeric-looking / percent-looking text that Excel would otherwise
# auto-convert to a Number cell, so those are written with a leading
# apostrophe (forces Text) and the cell style is reset to "Normal"
# afterwards so no stray quote-prefix/number-format style sticks around.
$changes = @'
D2|307.77
E2|-6.00%
D3|40.79
E3|-7.78%
D4|5.066
E4|-4.54%
D5|0.07839
E5|-6.20%
D6|4.322
E6|-1.85%
D7|1.664
E7|-14.72%
D8|0.9180
E8|-5.29%
E9|-6.55%
D10|0.1763
B11|BitrueCoin
C11|https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr
D11|0.04477
E11|-2.77%
B12|MandalaExchangeToken
C12|https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx
D12|0.08959
E12|-7.24%
D13|7.188
E13|-15.98%
D14|0.1058
E14|-0.09%
D15|0.001291
E15|0.12%
D16|0.005713
E16|-3.07%
D17|3.369
E17|-0.56%
E18|2.27%
D19|0.3367
E19|0.26%
D20|0.1386
E20|-0.20%
D21|0.2860
E21|11.00%
D22|0.04184
E22|0.77%
D23|0.001222
E23|-0.83%
D24|0.004134
E24|-6.23%
D25|0.0001230
E25|-5.50%
D26|0.0003002
E26|0.69%
D38|0.02411
E38|-10.92%
D39|0.05242
E39|-6.12%
D40|0.008011
E40|2.41%
D41|0.1350
E41|-4.25%
D42|0.007524
E42|3.00%
D43|0.001996
E43|-2.71%
D44|0.008080
E44|-7.03%
D45|0.3354
E45|-4.53%
D46|0.00006748
E46|-1.37%
D47|0.00000000756
E47|0.69%
D48|0.003423
E48|-1.86%
D49|0.004132
E49|16.95%
D50|0.00002116
E50|0.69%
D51|0.0002016
E51|0.69%
'@

$rows = $changes -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }

foreach ($row in $rows) {
    $parts = $row.Split('|', 2)
    $ref = $parts[0]
    $newValue = $parts[1]

    $cell = $ws.Range($ref)
    $col = $ref -replace '\d+$', ''

    if ($col -eq 'D' -or $col -eq 'E') {
        # Force text so "307.77" / "-6.00%" stay literal strings instead of
        # being reinterpreted as a Number/Percentage.
        $cell.Value = "'" + $newValue
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $newValue
    }
}
